$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("used")

# Delete row 2 of Sheet1 (the "used" id moves from the names pool), shifting rows up.
$ws1.Rows.Item(2).Delete()

# Append the newly-used entry to the "used" sheet.
$ws2.Cells.Item(43, 1).Value = "58jp6ymz"
$ws2.Cells.Item(43, 2).Value = "ChatGPT Image 2026年1月21日 21_22_39.png"
$ws2.Cells.Item(43, 3).Value = "2026-01-21 21:23:20"
